# Apply corrected data-cleaning edits for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): remove the bold/bordered "header" style and
#     clear the stray "Unnamed: 0" label that used to sit in A1 ---
$ws.Range("A1").Value = ""
$ws.Range("A1:BE1").Style = "Normal"

# --- Row 3: Revisit count ---
$ws.Range("I3").Value = 28
$ws.Range("J3").Value = 22
$ws.Range("Z3").Value = 26
$ws.Range("AA3").Value = 24
$ws.Range("AB3").Value = 25
$ws.Range("AP3").Value = 32

# --- Row 4: Fixation count ---
$ws.Range("I4").Value = 81
$ws.Range("J4").Value = 50
$ws.Range("Z4").Value = 75
$ws.Range("AA4").Value = 69
$ws.Range("AB4").Value = 62
$ws.Range("AP4").Value = 385

# --- Row 5: Dwell time (ms) ---
$ws.Range("I5").Value = 14867.03
$ws.Range("J5").Value = 8742.17
$ws.Range("Z5").Value = 13615.4
$ws.Range("AA5").Value = 12863.5
$ws.Range("AB5").Value = 11261.55
$ws.Range("AP5").Value = 83152.98

# --- Row 6: Dwell time (%) ---
$ws.Range("B6").Value = 0.15
$ws.Range("D6").Value = 0.11
$ws.Range("E6").Value = 0.51
$ws.Range("F6").Value = 0.16
$ws.Range("G6").Value = 0.13
$ws.Range("H6").Value = 0.06
$ws.Range("I6").Value = 5.74
$ws.Range("J6").Value = 3.37
$ws.Range("K6").Value = 0.63
$ws.Range("L6").Value = 0.77
$ws.Range("M6").Value = 0.18
$ws.Range("N6").Value = 0.25
$ws.Range("T6").Value = 0.26
$ws.Range("U6").Value = 0.53
$ws.Range("V6").Value = 0.15
$ws.Range("W6").Value = 0.11
$ws.Range("X6").Value = 0.26
$ws.Range("Z6").Value = 5.25
$ws.Range("AA6").Value = 4.96
$ws.Range("AB6").Value = 4.34
$ws.Range("AD6").Value = 0.15
$ws.Range("AE6").Value = 0.31
$ws.Range("AG6").Value = 0.07
$ws.Range("AJ6").Value = 0.21
$ws.Range("AK6").Value = 0.49
$ws.Range("AP6").Value = 32.08
$ws.Range("AQ6").Value = 0.23

# --- Row 7: Fixation duration (ms) ---
$ws.Range("I7").Value = 183.54
$ws.Range("J7").Value = 174.84
$ws.Range("Z7").Value = 181.54
$ws.Range("AA7").Value = 186.43
$ws.Range("AB7").Value = 181.64
$ws.Range("AP7").Value = 215.98

$wb.Save()
